$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (BTC)
$ws.Range("C3").Value = 896.05924946
$ws.Range("D3").Value = 123.63124776
$ws.Range("E3").Value = 772.4280017
$ws.Range("G3").Value = 896.0599999999999

# Row 4 (ATOM)
$ws.Range("F4").Value = 34312.7
$ws.Range("G4").Value = 528.47

# Row 5 (ALGO)
$ws.Range("C5").Value = 22.0802
$ws.Range("E5").Value = 21.843
$ws.Range("F5").Value = 13.3314
$ws.Range("G5").Value = 294.36

# Row 6 (ETH)
$ws.Range("F6").Value = 0.9287

# Row 7 (MATIC)
$ws.Range("F7").Value = 2325.29
